$d = $word.ActiveDocument

$replacements = @(
    @("2024-04-12 Friday", "2024-04-13 Saturday"),
    @("354÷4=", "287÷8="),
    @("445÷9=", "713÷9="),
    @("416÷6=", "644÷2="),
    @("862÷8=", "122÷3="),
    @("574÷8=", "387÷9="),
    @("501÷9=", "432÷7="),
    @("428÷4=", "464÷6="),
    @("927÷4=", "513÷3="),
    @("879÷8=", "156÷9="),
    @("825÷2=", "425÷4="),
    @("860÷3=", "684÷9="),
    @("969÷4=", "749÷2="),
    @("134÷3=", "228÷2="),
    @("501÷8=", "996÷9="),
    @("706÷5=", "371÷9="),
    @("144÷6=", "602÷3="),
    @("930÷9=", "522÷5="),
    @("142÷6=", "127÷2="),
    @("279÷2=", "744÷5="),
    @("869÷5=", "769÷5="),
    @("185÷7=", "120÷3="),
    @("275÷3=", "984÷2="),
    @("347÷6=", "560÷5="),
    @("225÷9=", "336÷2="),
    @("249÷8=", "582÷6=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
